$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315674185752869
$ws.Range("B1").Value = 2.104238271713257
$ws.Range("C1").Value = 4.760006904602051
$ws.Range("D1").Value = 3.516060829162598
$ws.Range("E1").Value = 1.3377525806427
